# Weekly refresh of the "Fruta, Vega Central Mapocho de Santiago - Coco" sheet.
# The underlying per-record fields (Fecha, Calidad, Volumen, Precio min/max/prom,
# Origen, Precio $/Kg) are reshuffled across the existing data rows (2-41) -
# i.e. each row ends up with the full record previously held by another row.
# Columns A:C (Mercado/Región), E:K (Codreg..Variedad), Q (Unidad) and T
# (Kg/unidad) are identical for every record, so re-writing the whole A:T
# block per row is equivalent to only touching the columns that actually
# differ, and is far less error-prone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 41
$firstCol = 1   # A
$lastCol  = 20  # T
$numCols  = $lastCol - $firstCol + 1
$numRows  = $lastRow - $firstRow + 1

# Snapshot every data row BEFORE any write, since several rows trade places
# with each other (the mapping below is a permutation, not independent
# per-row overwrites).
$snapshot = $ws.Range("A$firstRow`:T$lastRow").Value2

# target Excel row -> source Excel row (source row's A:T content is written
# into the target row)
$rowMap = @{
  2=39; 3=41; 4=4; 5=27; 6=12; 7=29; 8=2; 9=35; 10=24; 11=34;
  12=8; 13=19; 14=25; 15=28; 16=17; 17=37; 18=40; 19=33; 20=14;
  21=16; 22=21; 23=31; 24=15; 25=6; 26=9; 27=18; 28=10; 29=7;
  30=23; 31=3; 32=26; 33=30; 34=13; 35=22; 36=38; 37=11; 38=36;
  39=32; 40=5; 41=20
}

foreach ($targetRow in $firstRow..$lastRow) {
    $sourceRow = $rowMap[$targetRow]
    if ($sourceRow -eq $targetRow) {
        continue
    }
    $srcOffset = $sourceRow - $firstRow + 1

    $newRow = New-Object 'object[,]' 1,$numCols
    for ($c = 1; $c -le $numCols; $c++) {
        $newRow[0, $c-1] = $snapshot[$srcOffset, $c]
    }

    $ws.Range("A$targetRow`:T$targetRow").Value2 = $newRow
}
